$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$summary = $wb.Worksheets.Item(2)

# These "Number"/"PackageID" values are stored as text in the source file
# (see the numberStoredAsText ignored-error markup), so force Text format
# on each numeric-looking cell before writing it, to keep it as text
# rather than letting Excel auto-convert it to a number.
$ws.Range("F21").NumberFormat = "@"
$ws.Range("F21").Value = "10"

# New row 22
$ws.Range("C22").Value = "14_波浪浅紫洋桔梗_Wavy Light Purple Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g"
$ws.Range("F22").NumberFormat = "@"
$ws.Range("F22").Value = "10"

# New row 23
$ws.Range("C23").Value = "1_白洋桔梗_White Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g"
$ws.Range("F23").NumberFormat = "@"
$ws.Range("F23").Value = "10"

# New row 24
$ws.Range("A24").NumberFormat = "@"
$ws.Range("A24").Value = "6"
$ws.Range("C24").Value = "4_阳光粉洋桔梗_Sunshine Pink Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g"
$ws.Range("F24").NumberFormat = "@"
$ws.Range("F24").Value = "15"

# New row 25
$ws.Range("C25").Value = "602_康乃馨白_white_undefined_20stems"
$ws.Range("F25").NumberFormat = "@"
$ws.Range("F25").Value = "15"

# New row 26 (name contains an embedded line break)
$ws.Range("C26").Value = "509_翠珠粉_Didiscus caeruleus
pink_Trachymene Coerulea_1bunch"
$ws.Range("F26").NumberFormat = "@"
$ws.Range("F26").Value = "10"

# New row 27 (name contains an embedded line break)
$ws.Range("C27").Value = "510_翠珠白_Didiscus caeruleus 
white_Trachymene Coerulea_1bunch"
$ws.Range("F27").NumberFormat = "@"
$ws.Range("F27").Value = "20"

# New row 28
$ws.Range("C28").Value = "8_冰淇淋洋桔梗_Icecream Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g"
$ws.Range("F28").NumberFormat = "@"
$ws.Range("F28").Value = "10"

# New row 29
$ws.Range("C29").Value = "3_波浪白洋桔梗_Wavy White Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g"
$ws.Range("F29").NumberFormat = "@"
$ws.Range("F29").Value = "3"

# New row 30
$ws.Range("C30").Value = "616_康乃馨紫精灵_Purple Elves_undefined_20stems"
$ws.Range("F30").NumberFormat = "@"
$ws.Range("F30").Value = "5"

# New row 31 (no Number value)
$ws.Range("C31").Value = "238_苏菲宝贝_undefined_Rosa rugosa Thunb._10stems"

# Summary sheet: G2 encoded-number string grows with the new rows
$summary.Range("G2").NumberFormat = "@"
$summary.Range("G2").Value = "016111210101614101030101010305040105201010101515102010350"
